$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 285, shifting existing rows 285:400 down to 286:401
$ws.Rows.Item(285).Insert()

# Populate the newly inserted row 285 with the new record's data
$ws.Cells.Item(285, 1).Value = 5
$ws.Cells.Item(285, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(285, 3).Value = "Maule"
$ws.Cells.Item(285, 4).Value = 45009
$ws.Cells.Item(285, 5).Value = 7
$ws.Cells.Item(285, 6).Value = 100112009
$ws.Cells.Item(285, 7).Value = "Acelga"
$ws.Cells.Item(285, 8).Value = "Sin especificar"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 500
$ws.Cells.Item(285, 11).Value = 2500
$ws.Cells.Item(285, 12).Value = 2500
$ws.Cells.Item(285, 13).Value = 2500
$ws.Cells.Item(285, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(285, 15).Value = "Región del Maule"
$ws.Cells.Item(285, 16).Value = 625
$ws.Cells.Item(285, 17).Value = 4
$ws.Cells.Item(285, 18).Value = "Hortaliza"
